$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 23: date (A23) formatted like the other date cells (copy format from A22),
# and hours value (B23).
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A23").Value = 45172

$ws.Range("B23").Value = 6.5

$wb.Save()
